$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Locate the paragraph that currently ends with:
#      "Rewrite data from text file to JSON object to C# object"
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Rewrite data from text file to JSON object to C# object*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {

    # Append a new run "(DATA IS ALDREADY JSON?)" right after the
    # existing text, inside the same paragraph (before the paragraph
    # mark). Toggling Bold on/off for the freshly inserted text forces
    # Word to keep it as its own run instead of merging it back into
    # the preceding run (both runs end up with identical formatting,
    # matching the target markup).
    $insertPos = $target.Range.End - 1
    $insertRange = $d.Range($insertPos, $insertPos)
    $newText = "(DATA IS ALDREADY JSON?)"
    $insertRange.InsertAfter($newText)

    $newRunRange = $d.Range($insertPos, $insertPos + $newText.Length)
    $newRunRange.Bold = 1
    $newRunRange.Bold = 0

    # -------------------------------------------------------------
    # 2. Move the "_GoBack" bookmark so it ends up right after the
    #    text we just inserted (this mirrors what Word itself does:
    #    it always tracks the location of the most recent edit).
    # -------------------------------------------------------------
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }

    # A bookmark collapsed exactly on a paragraph-mark boundary can be
    # mis-anchored, so: insert a one-character placeholder after the
    # new text, build the bookmark around that (non-edge) character,
    # then delete the placeholder again. The bookmark collapses back
    # to a zero-length bookmark at the correct position.
    $markerPos = $target.Range.End - 1
    $markerRange = $d.Range($markerPos, $markerPos)
    $markerRange.InsertAfter("X")

    $bmRange = $d.Range($markerPos, $markerPos + 1)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $markerDelete = $d.Range($markerPos, $markerPos + 1)
    $markerDelete.Text = ""
}
